$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 16:59"

# Country name reassignments caused by the updated sort order
$ws.Range("A20").Value = "Italia"
$ws.Range("A21").Value = "Arabia Saudita"
$ws.Range("A25").Value = "Alemania"
$ws.Range("A26").Value = "Pakistan"
$ws.Range("A58").Value = "Moldavia"
$ws.Range("A59").Value = "Suiza"
$ws.Range("A60").Value = "Uzbekistan"

# Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes)
$ws.Range("B4").Value = 7839447
$ws.Range("C4").Value = 5684
$ws.Range("D4").Value = 5026389
$ws.Range("E4").Value = 2595229
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = 217829

$ws.Range("B17").Value = 477769
$ws.Range("C17").Value = 1753
$ws.Range("D17").Value = 450297
$ws.Range("E17").Value = 14252
$ws.Range("G17").Value = 53
$ws.Range("H17").Value = 13220

$ws.Range("B20").Value = 343770
$ws.Range("C20").Value = 5372
$ws.Range("D20").Value = 237549
$ws.Range("E20").Value = 70110
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 36111

$ws.Range("B21").Value = 338539
$ws.Range("C21").Value = 407
$ws.Range("D21").Value = 324282
$ws.Range("E21").Value = 9261
$ws.Range("G21").Value = 24
$ws.Range("H21").Value = 36083

$ws.Range("B25").Value = 318007
$ws.Range("C25").Value = 2493
$ws.Range("D25").Value = 269500
$ws.Range("E25").Value = 38833
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 9674

$ws.Range("B26").Value = 317595
$ws.Range("C26").Value = 661
$ws.Range("D26").Value = 302708
$ws.Range("E26").Value = 8335
$ws.Range("G26").Value = 8
$ws.Range("H26").Value = 6552

$ws.Range("D57").Value = 69854
$ws.Range("E57").Value = 4302
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 266

$ws.Range("B58").Value = 60833
$ws.Range("C58").Value = 918
$ws.Range("D58").Value = 43489
$ws.Range("E58").Value = 15902
$ws.Range("G58").Value = 18
$ws.Range("H58").Value = 1442

$ws.Range("B59").Value = 60368
$ws.Range("C59").Value = 1487
$ws.Range("D59").Value = 48400
$ws.Range("E59").Value = 9880
$ws.Range("H59").Value = 2088

$ws.Range("B60").Value = 60342
$ws.Range("C60").Value = 316
$ws.Range("D60").Value = 57221
$ws.Range("E60").Value = 2623
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 498

$ws.Range("B95").Value = 15294
$ws.Range("C95").Value = 73
$ws.Range("E95").Value = 3156

$ws.Range("B97").Value = 15066
$ws.Range("C97").Value = 167
$ws.Range("D97").Value = 9304
$ws.Range("E97").Value = 5349
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = 413

$ws.Range("B104").Value = 10835
$ws.Range("C104").Value = 13
$ws.Range("E104").Value = 317

$ws.Range("B117").Value = 7363
$ws.Range("C117").Value = 90
$ws.Range("D117").Value = 3102
$ws.Range("E117").Value = 4129
$ws.Range("G117").Value = 4
$ws.Range("H117").Value = 132

$ws.Range("B121").Value = 5943
$ws.Range("C121").Value = 26
$ws.Range("D121").Value = 5398
$ws.Range("E121").Value = 422

$ws.Range("B136").Value = 4504
$ws.Range("C136").Value = 16
$ws.Range("E136").Value = 1195

$ws.Range("B137").Value = 4491
$ws.Range("C137").Value = 106
$ws.Range("D137").Value = 3994
$ws.Range("E137").Value = 481

$ws.Range("B139").Value = 3989
$ws.Range("C139").Value = 97
$ws.Range("E139").Value = 982
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 43

$ws.Range("B143").Value = 3621
$ws.Range("C143").Value = 4
$ws.Range("D143").Value = 2489
$ws.Range("E143").Value = 1015

$ws.Range("B150").Value = 2696
$ws.Range("C150").Value = 128
$ws.Range("D150").Value = 1814
$ws.Range("E150").Value = 827
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 55

$ws.Range("B164").Value = 1800
$ws.Range("C164").Value = 14
$ws.Range("E164").Value = 832
$ws.Range("G164").Value = 1
$ws.Range("H164").Value = 42
